$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold font, border, centered) from H1 onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-20: I column = 1 (except row 20 = 9), J column mirrors H column (except row 20 = 9)
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}

$ws.Cells.Item(20, 9).Value = 9
$ws.Cells.Item(20, 10).Value = 9
